$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new "Po No" header in column G, matching the style of the existing
# header cells (e.g. F1).
$ws.Range("G1").Value = "Po No"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Move/record the active selection as shown in the target workbook.
$ws.Range("G6").Select()
